$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (shifts old D..H to F..J)
$ws.Columns("D:E").Insert()

# Update header row
$ws.Cells.Item(1, 4).Value = 'Best Score'
$ws.Cells.Item(1, 5).Value = 'Search Type'

# Update data rows 2-6
# Row 2
$ws.Cells.Item(2, 2).Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''Delta_PV'',
                                                   ''Delta_RT'', ''FullPath_HR'',
                                                   ''FullPath_V'', ''MT_HR'',
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LogisticRegression(class_weight=''balanced'',
                                    l1_ratio=0.6176871281126421, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Cells.Item(2, 3).Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__l1_ratio'': 0.6176871281126421, ''model__penalty'': ''elasticnet'', ''model__solver'': ''saga''}'
$ws.Cells.Item(2, 4).Value = 0.5526651281295245
$ws.Cells.Item(2, 5).Value = 'Tree-Parzen Estimator'
$ws.Cells.Item(2, 6).Value = 42
$ws.Cells.Item(2, 7).Value = 0.7121122738459804
$ws.Cells.Item(2, 8).Value = 0.4611594202898551
$ws.Cells.Item(2, 9).Value = '[1 0 1 0 0 1 1 1 1 1 1 1 1 0 1 0 0 0 1 0 1 1 0 0]'
$ws.Cells.Item(2, 10).Value = '[0 1 1 1 0 0 1 0 0 1 1 0 0 1 1 0 1 0 0 0 1 0 0 1]'
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(2).EntireRow.AutoFit()

# Row 3
$ws.Cells.Item(3, 2).Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''Delta_PV'',
                                                   ''Delta_RT'', ''FullPath_HR'',
                                                   ''FullPath_V'', ''MT_HR'',
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7f668c093970>),
                (''model'',
                 LogisticRegression(class_weight=''balanced'',
                                    l1_ratio=0.589923362819949, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Cells.Item(3, 3).Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7f67799fa220>, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__l1_ratio'': 0.589923362819949, ''model__penalty'': ''elasticnet'', ''model__solver'': ''saga''}'
$ws.Cells.Item(3, 4).Value = 0.5440036044929032
$ws.Cells.Item(3, 5).Value = 'Tree-Parzen Estimator'
$ws.Cells.Item(3, 6).Value = 69
$ws.Cells.Item(3, 7).Value = 0.6374807987711214
$ws.Cells.Item(3, 8).Value = 0.5833333333333334
$ws.Cells.Item(3, 9).Value = '[0 1 1 0 1 0 0 0 1 1 1 0 1 0 1 0 1 0 1 1 0 1 1 1]'
$ws.Cells.Item(3, 10).Value = '[1 1 1 1 0 1 0 1 0 1 0 0 1 0 1 0 1 1 1 0 0 0 1 1]'
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(3).EntireRow.AutoFit()

# Row 4
$ws.Cells.Item(4, 2).Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''Delta_PV'',
                                                   ''Delta_RT'', ''FullPath_HR'',
                                                   ''FullPath_V'', ''MT_HR'',
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LogisticRegression(class_weight=''balanced'',
                                    l1_ratio=0.09772327678985551, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Cells.Item(4, 3).Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__l1_ratio'': 0.09772327678985551, ''model__penalty'': ''elasticnet'', ''model__solver'': ''saga''}'
$ws.Cells.Item(4, 4).Value = 0.4847212296609927
$ws.Cells.Item(4, 5).Value = 'Tree-Parzen Estimator'
$ws.Cells.Item(4, 6).Value = 23
$ws.Cells.Item(4, 7).Value = 0.6887664964018915
$ws.Cells.Item(4, 8).Value = 0.6666666666666666
$ws.Cells.Item(4, 9).Value = '[0 0 1 0 0 1 0 1 1 1 1 1 1 1 1 0 0 0 1 0 1 1 1 0]'
$ws.Cells.Item(4, 10).Value = '[0 0 1 0 1 1 1 1 1 1 1 0 0 1 1 1 0 1 0 0 1 0 1 0]'
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(4).EntireRow.AutoFit()

# Row 5
$ws.Cells.Item(5, 2).Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''Delta_PV'',
                                                   ''Delta_RT'', ''FullPath_HR'',
                                                   ''FullPath_V'', ''MT_HR'',
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LogisticRegression(class_weight=''balanced'',
                                    l1_ratio=0.6897879621384697, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Cells.Item(5, 3).Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__l1_ratio'': 0.6897879621384697, ''model__penalty'': ''elasticnet'', ''model__solver'': ''saga''}'
$ws.Cells.Item(5, 4).Value = 0.5309397427926166
$ws.Cells.Item(5, 5).Value = 'Tree-Parzen Estimator'
$ws.Cells.Item(5, 6).Value = 99
$ws.Cells.Item(5, 7).Value = 0.7007709474538447
$ws.Cells.Item(5, 8).Value = 0.5034965034965034
$ws.Cells.Item(5, 9).Value = '[0 0 1 1 0 1 1 1 1 0 1 1 0 1 1 0 0 1 0 1 1 0 0 1]'
$ws.Cells.Item(5, 10).Value = '[0 0 1 0 1 1 1 1 0 1 0 1 1 0 1 0 0 0 1 0 1 1 0 0]'
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(5).EntireRow.AutoFit()

# Row 6
$ws.Cells.Item(6, 2).Value = 'Pipeline(steps=[(''scaler'',
                 ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                                   transformers=[(''StandardScaler'',
                                                  StandardScaler(),
                                                  [''AE_HR'', ''AE_V'',
                                                   ''AbsOffAxis_HR'',
                                                   ''AbsOffAxis_V'',
                                                   ''AbsOnAxis_HR'',
                                                   ''AbsOnAxis_V'', ''BallPath_HR'',
                                                   ''BallPath_V'', ''CMT_HR'',
                                                   ''CMT_V'', ''Corrective_HR'',
                                                   ''Corrective_V'', ''Delta_AE'',
                                                   ''Delta_Fullpath'', ''Delta_MT'',
                                                   ''Delta_OffAxis'',
                                                   ''Delta_OnAxis'', ''Delta_PV'',
                                                   ''Delta_RT'', ''FullPath_HR'',
                                                   ''FullPath_V'', ''MT_HR'',
                                                   ''MT_V'', ''PeakV_HR'',
                                                   ''PeakV_V'', ''RT_HR'', ''RT_V'',
                                                   ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])])),
                (''selector'', None),
                (''model'',
                 LogisticRegression(class_weight=''balanced'',
                                    l1_ratio=0.9890779091209877, max_iter=1000,
                                    penalty=''elasticnet'', random_state=42,
                                    solver=''saga''))])'
$ws.Cells.Item(6, 3).Value = '{''selector'': None, ''scaler'': ColumnTransformer(n_jobs=-1, remainder=''passthrough'',
                  transformers=[(''StandardScaler'', StandardScaler(),
                                 [''AE_HR'', ''AE_V'', ''AbsOffAxis_HR'',
                                  ''AbsOffAxis_V'', ''AbsOnAxis_HR'', ''AbsOnAxis_V'',
                                  ''BallPath_HR'', ''BallPath_V'', ''CMT_HR'',
                                  ''CMT_V'', ''Corrective_HR'', ''Corrective_V'',
                                  ''Delta_AE'', ''Delta_Fullpath'', ''Delta_MT'',
                                  ''Delta_OffAxis'', ''Delta_OnAxis'', ''Delta_PV'',
                                  ''Delta_RT'', ''FullPath_HR'', ''FullPath_V'',
                                  ''MT_HR'', ''MT_V'', ''PeakV_HR'', ''PeakV_V'',
                                  ''RT_HR'', ''RT_V'', ''TMT_HR'', ''TMT_V'', ''VE_HR'', ...])]), ''model__l1_ratio'': 0.9890779091209877, ''model__penalty'': ''elasticnet'', ''model__solver'': ''saga''}'
$ws.Cells.Item(6, 4).Value = 0.5562447465637667
$ws.Cells.Item(6, 5).Value = 'Tree-Parzen Estimator'
$ws.Cells.Item(6, 6).Value = 89
$ws.Cells.Item(6, 7).Value = 0.7007709474538447
$ws.Cells.Item(6, 8).Value = 0.5440579710144928
$ws.Cells.Item(6, 9).Value = '[1 0 1 0 1 1 0 0 1 1 0 1 0 1 1 1 1 1 0 1 0 0 1 0]'
$ws.Cells.Item(6, 10).Value = '[0 0 1 1 1 0 0 0 1 0 1 1 0 1 1 0 0 0 0 1 0 1 0 1]'
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(6).EntireRow.AutoFit()
